# Update countries & provincias Spain
# - New Covid-19 figures for several countries (US, Germany, Brazil, Israel, Tunisia, Trinidad & Tobago)
# - Fresh data for Granada, which is reinserted into the ranked list right after Belice,
#   pushing Fiyi -> San Cristobal y Nieves down by one row
# - "Datos actualizados" timestamp bumped from 22:22 to 22:52

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 25 de Abril de 2020 a las 22:52"

$cols = @("B","C","D","E","F","G","H")

$rows = @(
    @{ Row=4; Name="Estados Unidos"; Vals=@(953851, 28619, 116015, 783985, 15110, 1658, 53851) },
    @{ Row=8; Name="Alemania"; Vals=@(156126, 1127, 109800, 40480, 2908, 86, 5846) },
    @{ Row=14; Name="Brasil"; Vals=@(58616, 5621, 27655, 26945, 8318, 346, 4016) },
    @{ Row=26; Name="Israel"; Vals=@(15298, 240, 6435, 8664, 127, 5, 199) },
    @{ Row=90; Name="Tunez"; Vals=@(939, 17, 207, 694, 19, 0, 38) },
    @{ Row=144; Name="Trinidad yTobago"; Vals=@(115, 0, 53, 54, 0, 0, 8) },
    @{ Row=184; Name="Granada"; Vals=@(18, 3, 7, 11, 4, 0, 0) },
    @{ Row=185; Name="Fiyi"; Vals=@(18, 0, 10, 8, 0, 0, 0) },
    @{ Row=186; Name="Nueva Caledonia"; Vals=@(18, 0, 17, 1, 1, 0, 0) },
    @{ Row=187; Name="Islas Virgenes de los Estados Unidos"; Vals=@(17, 0, 0, 17, 0, 0, 0) },
    @{ Row=188; Name="Namibia"; Vals=@(16, 0, 7, 9, 0, 0, 0) },
    @{ Row=190; Name="Republica de Africa Central"; Vals=@(16, 0, 10, 6, 0, 0, 0) },
    @{ Row=191; Name="Curazao"; Vals=@(16, 0, 11, 4, 0, 0, 1) },
    @{ Row=192; Name="San Cristobal y Nieves"; Vals=@(15, 0, 2, 13, 0, 0, 0) }
)

foreach ($r in $rows) {
    $ws.Range("A" + $r.Row).Value = $r.Name
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $r.Row).Value = $r.Vals[$i]
    }
}
